{"js": "// Insert a new paragraph \"\u554a\u5927\u5927\" right after the first paragraph (\"Aaabbbb\"),\n// then change the text of what was the second paragraph (also \"\u554a\u5927\u5927\")\n// to \"Dev\", leaving its bookmark intact.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The first paragraph (\"Aaabbbb\") gets a new sibling paragraph right after\n// it. Anchoring the insert on the first paragraph (rather than inserting\n// \"Before\" the second paragraph) makes the new paragraph mark inherit the\n// east-Asian formatting hint from paragraph 1, matching the target edit.\nconst firstParagraph = paragraphs.items[0];\nconst insertedParagraph = firstParagraph.insertParagraph(\"\u554a\u5927\u5927\", Word.InsertLocation.after);\n\n// Re-load the paragraph collection so we can reliably reach the paragraph\n// that originally held \"\u554a\u5927\u5927\" (now the third paragraph) and flip its text\n// to \"Dev\" while keeping its own run/paragraph formatting and the\n// _GoBack bookmark untouched.\nconst refreshedParagraphs = body.paragraphs;\nrefreshedParagraphs.load(\"items\");\nawait context.sync();\n\nconst targetParagraph = refreshedParagraphs.items[2];\nconst targetRange = targetParagraph.getRange(Word.RangeLocation.whole);\ntargetRange.insertText(\"Dev\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Insert a new paragraph \"\u554a\u5927\u5927\" right after the first paragraph (\"Aaabbbb\"),\n# then change the text of what was the second paragraph (itself \"\u554a\u5927\u5927\")\n# to \"Dev\", leaving its bookmark (_GoBack) untouched.\n\n$d = $word.ActiveDocument\n\n# Anchoring the new paragraph mark on paragraph 1 (rather than inserting\n# before paragraph 2) makes it inherit paragraph 1's east-Asian formatting\n# hint, matching the target edit.\n$firstParagraph = $d.Paragraphs.Item(1)\n$firstParagraph.Range.InsertParagraphAfter()\n\n# The document now has three paragraphs: \"Aaabbbb\", a fresh empty one, and\n# the original \"\u554a\u5927\u5927\" paragraph (with its bookmark) pushed down to #3.\n$newParagraph = $d.Paragraphs.Item(2)\n$newParagraph.Range.Text = \"\u554a\u5927\u5927\"\n\n# Paragraph 3 is still the original paragraph (same bookmark intact);\n# just swap its run text from \"\u554a\u5927\u5927\" to \"Dev\".\n$originalParagraph = $d.Paragraphs.Item(3)\n$originalParagraph.Range.Text = \"Dev\"\n"}
